$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 891, shifting existing rows 891:989 down to 892:990
$ws.Rows("891:891").Insert()

# Populate the newly inserted row 891 with the new weekly data point
$ws.Range("A891").Value = 4
$ws.Range("B891").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C891").Value = "Los Lagos"
$ws.Range("D891").Value = 45194
$ws.Range("E891").Value = 10
$ws.Range("F891").Value = 100112033
$ws.Range("G891").Value = "Lechuga"
$ws.Range("H891").Value = "Escarola"
$ws.Range("I891").Value = "Primera"
$ws.Range("J891").Value = 300
$ws.Range("K891").Value = 17000
$ws.Range("L891").Value = 17000
$ws.Range("M891").Value = 17000
$ws.Range("N891").Value = "`$/caja 15 unidades"
$ws.Range("O891").Value = "Región de Coquimbo"
$ws.Range("P891").Value = 1133
$ws.Range("Q891").Value = 15
$ws.Range("R891").Value = "Hortaliza"
